# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2023-10-21 Saturday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2023-10-22 Sunday", 2)

# Update the practice-problem table. Each cell is addressed positionally
# (row, column) via the Tables object model so that the update is
# unambiguous regardless of whether a new value happens to collide with
# another cell's old/new text elsewhere in the table.
$tbl = $d.Tables.Item(1)

$newValues = @{
    1  = @("85÷4=21, 1", "53÷7=7, 4", "44÷6=7, 2", "73÷6=12, 1", "89÷9=9, 8")
    5  = @("86÷7=12, 2", "58÷8=7, 2", "31÷4=7, 3", "57÷3=19, 0", "84÷6=14, 0")
    9  = @("99÷5=19, 4", "93÷3=31, 0", "93÷9=10, 3", "98÷2=49, 0", "40÷8=5, 0")
    13 = @("34÷8=4, 2", "89÷3=29, 2", "54÷7=7, 5", "37÷8=4, 5", "33÷7=4, 5")
    17 = @("72÷8=9, 0", "14÷6=2, 2", "72÷3=24, 0", "43÷6=7, 1", "51÷9=5, 6")
}

foreach ($rowIndex in $newValues.Keys) {
    $rowValues = $newValues[$rowIndex]
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $cell = $tbl.Cell($rowIndex, $col)
        $cell.Range.Text = $rowValues[$col - 1]
    }
}
